# Apply the updated coin price / 1h-volume figures from the latest scrape.
# Values are stored as literal text (leading apostrophe) to match the
# workbook's existing string-typed Price/Volume columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'34.492.75"
$ws.Range("D3").Value = "'1.810.56"
$ws.Range("E3").Value = "'  +0.34%  "
$ws.Range("E4").Value = "'  -0.25%  "
$ws.Range("D5").Value = "'225.66"
$ws.Range("E5").Value = "'  -0.91%  "
$ws.Range("D6").Value = "'0.594"
$ws.Range("E6").Value = "'  +3.04%  "
$ws.Range("D8").Value = "'38.38"
$ws.Range("E8").Value = "'  +6.58%  "
$ws.Range("E9").Value = "'  -3.90%  "
$ws.Range("D10").Value = "'0.0675"
$ws.Range("E10").Value = "'  -2.51%  "
$ws.Range("D11").Value = "'0.0973"
$ws.Range("E11").Value = "'  +0.57%  "
$ws.Range("D12").Value = "'2.072.87"
$ws.Range("E12").Value = "'  +0.43%  "
$ws.Range("D13").Value = "'11.22"
$ws.Range("E13").Value = "'  -3.69%  "
$ws.Range("D14").Value = "'1.813.53"
$ws.Range("E14").Value = "'  +0.38%  "
$ws.Range("D15").Value = "'0.633"
$ws.Range("E15").Value = "'  -1.71%  "
$ws.Range("D16").Value = "'34.474.51"
$ws.Range("E16").Value = "'  +0.46%  "
$ws.Range("E17").Value = "'  -1.34%  "
$ws.Range("D18").Value = "'68.33"
$ws.Range("D19").Value = "'243.39"
$ws.Range("E19").Value = "'  -0.72%  "
$ws.Range("D20").Value = "'0.0₃0773"
$ws.Range("E20").Value = "'  -2.64%  "
$ws.Range("D21").Value = "'11.21"
$ws.Range("E21").Value = "'  -2.46%  "
$ws.Range("E22").Value = "'  -0.16%  "
$ws.Range("E23").Value = "'  -1.33%  "
$ws.Range("E24").Value = "'  +3.71%  "
$ws.Range("D25").Value = "'170.74"
$ws.Range("E25").Value = "'  -1.15%  "
$ws.Range("E26").Value = "'  -1.94%  "
$ws.Range("D27").Value = "'17.67"
$ws.Range("E27").Value = "'  +4.85%  "
$ws.Range("E28").Value = "'  +2.24%  "
$ws.Range("E29").Value = "'  -0.22%  "
$ws.Range("E30").Value = "'  -1.00%  "
$ws.Range("E32").Value = "'  -2.67%  "
$ws.Range("D33").Value = "'3.86"
$ws.Range("E33").Value = "'  -4.38%  "
$ws.Range("E34").Value = "'  +0.15%  "
$ws.Range("D35").Value = "'1.356.45"
$ws.Range("E35").Value = "'  -2.67%  "
$ws.Range("E36").Value = "'  -4.49%  "
$ws.Range("E37").Value = "'  -0.73%  "
$ws.Range("D38").Value = "'0.0188"
$ws.Range("E38").Value = "'  -1.10%  "
$ws.Range("E39").Value = "'  -4.87%  "
$ws.Range("D40").Value = "'2.45"
$ws.Range("E40").Value = "'  +1.26%  "
$ws.Range("D41").Value = "'0.955"
$ws.Range("E41").Value = "'  -0.60%  "
$ws.Range("D42").Value = "'82.03"
$ws.Range("E42").Value = "'  +0.48%  "
$ws.Range("D43").Value = "'1.22"
$ws.Range("E43").Value = "'  -1.40%  "
$ws.Range("E44").Value = "'  -1.00%  "
$ws.Range("D45").Value = "'13.73"
$ws.Range("E45").Value = "'  +1.36%  "
$ws.Range("E46").Value = "'  +1.57%  "
$ws.Range("D47").Value = "'1.973.92"
$ws.Range("E47").Value = "'  +0.45%  "
$ws.Range("E48").Value = "'  -3.96%  "
$ws.Range("E49").Value = "'  -0.20%  "
$ws.Range("D50").Value = "'102.64"
$ws.Range("E50").Value = "'  -2.13%  "
$ws.Range("D51").Value = "'0.0₆0122"
$ws.Range("E51").Value = "'  -4.74%  "
